$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Drop the trailing "Bibliografia" row (old row 22) completely -
#    it is the last row, so deleting it simply shrinks the sheet's
#    used range from C22 down to C21, matching the new <dimension>.
# ------------------------------------------------------------------
$ws.Rows.Item(22).Delete() | Out-Null

# ------------------------------------------------------------------
# 2) Row 10: the long Portuguese "Objetivos" paragraph in B/C is
#    replaced by the docente text.
# ------------------------------------------------------------------
$ws.Range("B10").Value = "9149242 - Fernando Catalani"
$ws.Range("C10").Value = "9149242 - Fernando Catalani"

# ------------------------------------------------------------------
# 3) Row 13 becomes "Programa resumido:" / "Semestral" (A13 is new).
# ------------------------------------------------------------------
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# ------------------------------------------------------------------
# 4) Row 14 keeps only the "Short syllabus:" label - drop B14/C14
#    entirely (Clear, not just ClearContents, so the <c> node itself
#    disappears instead of leaving an empty styled cell behind).
# ------------------------------------------------------------------
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Clear() | Out-Null
$ws.Range("C14").Clear() | Out-Null

# ------------------------------------------------------------------
# 5) Row 15 becomes "Programa:" / "01/01/2018" (B15/C15 are new
#    cells - copy the date-like text+format from B8/C8 so it keeps
#    its original text data type instead of being re-interpreted as
#    a date, and lands with the correct column style).
# ------------------------------------------------------------------
$ws.Range("A15").Value = "Programa:"

$ws.Range("B8").Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4163) | Out-Null
$ws.Range("B8").Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4122) | Out-Null

$ws.Range("C8").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4163) | Out-Null
$ws.Range("C8").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 6) Row 16 becomes "Syllabus:" with the English syllabus text.
# ------------------------------------------------------------------
$ws.Range("A16").Value = "Syllabus:"
$englishSyllabus = "1) Simple measures. Error Estimation of measures. Error propagation and significant figures.2) Construction of Tables and Graphs. Linearization.3) Introduction to the method of squares linear regression minimum.4) Kinematics. Rectilinear motion and uniformly varied motion. Free fall.5) Statics. Equilibrium of a material point. 6) Friction.7) Hooke's Law. Young´s Modulus.8) Energy conservation. Conservation Concept of Energy Mechanics. Mass-spring system.9) Shocks."
$ws.Range("B16").Value = $englishSyllabus
$ws.Range("C16").Value = $englishSyllabus

# ------------------------------------------------------------------
# 7) Row 17 keeps only the "Avaliação:" label - drop B17/C17.
# ------------------------------------------------------------------
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").Clear() | Out-Null
$ws.Range("C17").Clear() | Out-Null

# ------------------------------------------------------------------
# 8) Row 18 becomes "Método:" / "9149242 - Fernando Catalani"
#    (B18/C18 are new cells - copy column formatting so B18 doesn't
#    fall back to the wrong default style).
# ------------------------------------------------------------------
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "9149242 - Fernando Catalani"
$ws.Range("C18").Value = "9149242 - Fernando Catalani"

$ws.Range("B10").Copy() | Out-Null
$ws.Range("B18").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 9) Remaining label shifts (values in B/C already correct).
# ------------------------------------------------------------------
$ws.Range("A19").Value = "Critério:"
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("A21").Value = "Bibliografia:"

# ------------------------------------------------------------------
# 10) Row heights.
# ------------------------------------------------------------------
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(17).AutoFit() | Out-Null
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
